# Auto-generated Excel COM-interop script applying the diff
$wb = $excel.ActiveWorkbook

# --- 展览: simple numeric bumps ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 201
$ws1.Cells.Item(4, 6).Value = 69
$ws1.Cells.Item(7, 6).Value = 72
$ws1.Cells.Item(13, 6).Value = 2297
$ws1.Cells.Item(16, 6).Value = 519
$ws1.Cells.Item(17, 6).Value = 530
$ws1.Cells.Item(20, 6).Value = 42
$ws1.Cells.Item(22, 6).Value = 1775
$ws1.Cells.Item(23, 6).Value = 3909
$ws1.Cells.Item(32, 6).Value = 97
$ws1.Cells.Item(36, 6).Value = 681
$ws1.Cells.Item(38, 6).Value = 405

# --- 展览: rows 25-31 content refresh (new/updated conventions) ---
# row 25
$ws1.Cells.Item(25, 3).Value = "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会"
$ws1.Cells.Item(25, 4).Value = "兴国路恒大帝景西门 江西长庚控股有限公司"
$ws1.Cells.Item(25, 5).Value = "2024.07.28 11:00-07.28 17:00"
$ws1.Cells.Item(25, 6).Value = 58
$ws1.Cells.Item(25, 7).Value = 56
$ws1.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85688"
$ws1.Cells.Item(25, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png"
# row 26
$ws1.Cells.Item(26, 2).NumberFormat = "@"
$ws1.Cells.Item(26, 2).Value = "2024-08-03"
$ws1.Cells.Item(26, 3).Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$ws1.Cells.Item(26, 4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws1.Cells.Item(26, 5).Value = "2024.08.03 09:00-08.04 17:30"
$ws1.Cells.Item(26, 6).Value = 1165
$ws1.Cells.Item(26, 7).Value = 64
$ws1.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$ws1.Cells.Item(26, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"
# row 27
$ws1.Cells.Item(27, 3).Value = "吉安·COMIC LIFE周年庆典"
$ws1.Cells.Item(27, 4).Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws1.Cells.Item(27, 5).Value = "2024.08.03 09:30-08.03 18:00"
$ws1.Cells.Item(27, 6).Value = 221
$ws1.Cells.Item(27, 7).Value = 46.6
$ws1.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87164"
$ws1.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg"
# row 28
$ws1.Cells.Item(28, 3).Value = "景德镇·第十五届瓷都ACG动漫游戏博览会"
$ws1.Cells.Item(28, 4).Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws1.Cells.Item(28, 5).Value = "2024.08.03 09:00-08.04 17:00"
$ws1.Cells.Item(28, 6).Value = 2067
$ws1.Cells.Item(28, 7).Value = 55
$ws1.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86341"
$ws1.Cells.Item(28, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png"
# row 29
$ws1.Cells.Item(29, 3).Value = "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票"
$ws1.Cells.Item(29, 5).Value = "2024.08.03 08:30-08.03 17:00"
$ws1.Cells.Item(29, 6).Value = 568
$ws1.Cells.Item(29, 7).Value = "已售罄"
$ws1.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85981"
$ws1.Cells.Item(29, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png"
# row 30
$ws1.Cells.Item(30, 3).Value = "樟树·第二届静卿国风动漫文化展览会"
$ws1.Cells.Item(30, 4).Value = "杏佛路89号 樟树银河国际酒店"
$ws1.Cells.Item(30, 5).Value = "2024.08.03 09:00-08.03 17:00"
$ws1.Cells.Item(30, 6).Value = 465
$ws1.Cells.Item(30, 7).Value = 45
$ws1.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86683"
$ws1.Cells.Item(30, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg"
# row 31
$ws1.Cells.Item(31, 3).Value = "萍乡·AU9夏至国漫展"
$ws1.Cells.Item(31, 4).Value = "金陵东路18号 萍乡市体育馆"
$ws1.Cells.Item(31, 5).Value = "2024.08.03 10:00-08.03 17:00"
$ws1.Cells.Item(31, 6).Value = 63
$ws1.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86453"
$ws1.Cells.Item(31, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg"

# --- 演出: simple numeric bumps ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 26
$ws2.Cells.Item(2, 7).Value = 80

# --- 全部类型: simple numeric bumps ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 201
$ws4.Cells.Item(4, 6).Value = 69
$ws4.Cells.Item(7, 6).Value = 72
$ws4.Cells.Item(13, 6).Value = 2297
$ws4.Cells.Item(15, 6).Value = 26
$ws4.Cells.Item(15, 7).Value = 80
$ws4.Cells.Item(17, 6).Value = 519
$ws4.Cells.Item(18, 6).Value = 530
$ws4.Cells.Item(21, 6).Value = 42
$ws4.Cells.Item(23, 6).Value = 1775
$ws4.Cells.Item(24, 6).Value = 3909
$ws4.Cells.Item(33, 6).Value = 97
$ws4.Cells.Item(37, 6).Value = 681
$ws4.Cells.Item(39, 6).Value = 405

# --- 全部类型: rows 26-32 content refresh (mirrors 展览 25-31, offset +1) ---
# row 26
$ws4.Cells.Item(26, 3).Value = "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会"
$ws4.Cells.Item(26, 4).Value = "兴国路恒大帝景西门 江西长庚控股有限公司"
$ws4.Cells.Item(26, 5).Value = "2024.07.28 11:00-07.28 17:00"
$ws4.Cells.Item(26, 6).Value = 58
$ws4.Cells.Item(26, 7).Value = 56
$ws4.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85688"
$ws4.Cells.Item(26, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png"
# row 27
$ws4.Cells.Item(27, 2).NumberFormat = "@"
$ws4.Cells.Item(27, 2).Value = "2024-08-03"
$ws4.Cells.Item(27, 3).Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$ws4.Cells.Item(27, 4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws4.Cells.Item(27, 5).Value = "2024.08.03 09:00-08.04 17:30"
$ws4.Cells.Item(27, 6).Value = 1165
$ws4.Cells.Item(27, 7).Value = 64
$ws4.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$ws4.Cells.Item(27, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"
# row 28
$ws4.Cells.Item(28, 3).Value = "吉安·COMIC LIFE周年庆典"
$ws4.Cells.Item(28, 4).Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws4.Cells.Item(28, 5).Value = "2024.08.03 09:30-08.03 18:00"
$ws4.Cells.Item(28, 6).Value = 221
$ws4.Cells.Item(28, 7).Value = 46.6
$ws4.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87164"
$ws4.Cells.Item(28, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg"
# row 29
$ws4.Cells.Item(29, 3).Value = "景德镇·第十五届瓷都ACG动漫游戏博览会"
$ws4.Cells.Item(29, 4).Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws4.Cells.Item(29, 5).Value = "2024.08.03 09:00-08.04 17:00"
$ws4.Cells.Item(29, 6).Value = 2067
$ws4.Cells.Item(29, 7).Value = 55
$ws4.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86341"
$ws4.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png"
# row 30
$ws4.Cells.Item(30, 3).Value = "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票"
$ws4.Cells.Item(30, 5).Value = "2024.08.03 08:30-08.03 17:00"
$ws4.Cells.Item(30, 6).Value = 568
$ws4.Cells.Item(30, 7).Value = "已售罄"
$ws4.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85981"
$ws4.Cells.Item(30, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png"
# row 31
$ws4.Cells.Item(31, 3).Value = "樟树·第二届静卿国风动漫文化展览会"
$ws4.Cells.Item(31, 4).Value = "杏佛路89号 樟树银河国际酒店"
$ws4.Cells.Item(31, 5).Value = "2024.08.03 09:00-08.03 17:00"
$ws4.Cells.Item(31, 6).Value = 465
$ws4.Cells.Item(31, 7).Value = 45
$ws4.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86683"
$ws4.Cells.Item(31, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg"
# row 32
$ws4.Cells.Item(32, 3).Value = "萍乡·AU9夏至国漫展"
$ws4.Cells.Item(32, 4).Value = "金陵东路18号 萍乡市体育馆"
$ws4.Cells.Item(32, 5).Value = "2024.08.03 10:00-08.03 17:00"
$ws4.Cells.Item(32, 6).Value = 63
$ws4.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86453"
$ws4.Cells.Item(32, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg"
